$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148, shifting existing rows 148:223 down to 149:224
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new weekly record
$ws.Cells.Item(148, 1).Value = 11
$ws.Cells.Item(148, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(148, 3).Value = "Bíobío"
$ws.Cells.Item(148, 4).Value = 44875
$ws.Cells.Item(148, 5).Value = 8
$ws.Cells.Item(148, 6).Value = 100112003
$ws.Cells.Item(148, 7).Value = "Ajo"
$ws.Cells.Item(148, 8).Value = "Chino"
$ws.Cells.Item(148, 9).Value = "1a (guarda)"
$ws.Cells.Item(148, 10).Value = 250
$ws.Cells.Item(148, 11).Value = 13000
$ws.Cells.Item(148, 12).Value = 14000
$ws.Cells.Item(148, 13).Value = 13520
$ws.Cells.Item(148, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(148, 15).Value = "China"
$ws.Cells.Item(148, 16).Value = 1352
$ws.Cells.Item(148, 17).Value = 10
$ws.Cells.Item(148, 18).Value = "Hortaliza"
